$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated odds values per row (columns G..AJ = 7..36).
# Row 2
$ws.Cells.Item(2, 16).Value = 1.5
$ws.Cells.Item(2, 17).Value = 2.37
$ws.Cells.Item(2, 19).Value = 1.62

# Row 3
$ws.Cells.Item(3, 16).Value = 1.54
$ws.Cells.Item(3, 19).Value = 1.67

# Row 4
$ws.Cells.Item(4, 7).Value = 3.4
$ws.Cells.Item(4, 8).Value = 3
$ws.Cells.Item(4, 9).Value = 2.3
$ws.Cells.Item(4, 10).Value = 1.11
$ws.Cells.Item(4, 11).Value = 6.5
$ws.Cells.Item(4, 12).Value = 1.44
$ws.Cells.Item(4, 13).Value = 2.63
$ws.Cells.Item(4, 14).Value = 2.5
$ws.Cells.Item(4, 15).Value = 1.5
$ws.Cells.Item(4, 16).Value = 1.54
$ws.Cells.Item(4, 17).Value = 2.25
$ws.Cells.Item(4, 19).Value = 1.67
$ws.Cells.Item(4, 20).Value = 8
$ws.Cells.Item(4, 22).Value = 13
$ws.Cells.Item(4, 23).Value = 41
$ws.Cells.Item(4, 24).Value = 34
$ws.Cells.Item(4, 26).Value = 6.5
$ws.Cells.Item(4, 27).Value = 6
$ws.Cells.Item(4, 28).Value = 17

# Row 5
$ws.Cells.Item(5, 19).Value = 1.67

# Row 6
$ws.Cells.Item(6, 7).Value = 3.3
$ws.Cells.Item(6, 8).Value = 3.3
$ws.Cells.Item(6, 9).Value = 2.15
$ws.Cells.Item(6, 18).Value = 1.83
$ws.Cells.Item(6, 19).Value = 1.83
$ws.Cells.Item(6, 20).Value = 9.5
$ws.Cells.Item(6, 22).Value = 12
$ws.Cells.Item(6, 23).Value = 34
$ws.Cells.Item(6, 30).Value = 7.5
$ws.Cells.Item(6, 31).Value = 10
$ws.Cells.Item(6, 33).Value = 19
$ws.Cells.Item(6, 34).Value = 19
$ws.Cells.Item(6, 36).Value = 251

# Row 9
$ws.Cells.Item(9, 7).Value = 1.7
$ws.Cells.Item(9, 8).Value = 3.6
$ws.Cells.Item(9, 10).Value = 1.06
$ws.Cells.Item(9, 11).Value = 9.5
$ws.Cells.Item(9, 15).Value = 1.8
$ws.Cells.Item(9, 26).Value = 9.5
$ws.Cells.Item(9, 30).Value = 13

# Row 12
$ws.Cells.Item(12, 7).Value = 2.92
$ws.Cells.Item(12, 8).Value = 3.2
$ws.Cells.Item(12, 16).Value = 1.38
$ws.Cells.Item(12, 17).Value = 2.47
$ws.Cells.Item(12, 18).Value = 1.74
$ws.Cells.Item(12, 19).Value = 1.98
$ws.Cells.Item(12, 20).Value = 8
$ws.Cells.Item(12, 22).Value = 9
$ws.Cells.Item(12, 23).Value = 28
$ws.Cells.Item(12, 24).Value = 20
$ws.Cells.Item(12, 25).Value = 25
$ws.Cells.Item(12, 28).Value = 11.25
$ws.Cells.Item(12, 30).Value = 6.7

# Row 13
$ws.Cells.Item(13, 7).Value = 2.02
$ws.Cells.Item(13, 8).Value = 3.45
$ws.Cells.Item(13, 9).Value = 3.05
$ws.Cells.Item(13, 12).Value = 1.24
$ws.Cells.Item(13, 13).Value = 3.65
$ws.Cells.Item(13, 14).Value = 1.78
$ws.Cells.Item(13, 15).Value = 1.82
$ws.Cells.Item(13, 18).Value = 1.67
$ws.Cells.Item(13, 19).Value = 2.08
$ws.Cells.Item(13, 20).Value = 6.7
$ws.Cells.Item(13, 22).Value = 7.4
$ws.Cells.Item(13, 23).Value = 15
$ws.Cells.Item(13, 24).Value = 13
$ws.Cells.Item(13, 25).Value = 21
$ws.Cells.Item(13, 26).Value = 10.5
$ws.Cells.Item(13, 27).Value = 5.9
$ws.Cells.Item(13, 28).Value = 11.75
$ws.Cells.Item(13, 29).Value = 45
$ws.Cells.Item(13, 30).Value = 8.5
$ws.Cells.Item(13, 31).Value = 13.5
$ws.Cells.Item(13, 32).Value = 9.25
$ws.Cells.Item(13, 33).Value = 30
$ws.Cells.Item(13, 35).Value = 26
$ws.Cells.Item(13, 36).Value = 300

# Row 14
$ws.Cells.Item(14, 7).Value = 1.02
$ws.Cells.Item(14, 8).Value = 10.25
$ws.Cells.Item(14, 9).Value = 40
$ws.Cells.Item(14, 20).Value = 19
$ws.Cells.Item(14, 21).Value = 9.25
$ws.Cells.Item(14, 22).Value = 23
$ws.Cells.Item(14, 23).Value = 6.3
$ws.Cells.Item(14, 24).Value = 15.5
$ws.Cells.Item(14, 25).Value = 70
$ws.Cells.Item(14, 26).Value = 32
$ws.Cells.Item(14, 27).Value = 50
$ws.Cells.Item(14, 28).Value = 110
$ws.Cells.Item(14, 29).Value = 450
$ws.Cells.Item(14, 30).Value = 350
$ws.Cells.Item(14, 32).Value = 400
$ws.Cells.Item(14, 35).Value = 450

# Row 15
$ws.Cells.Item(15, 16).Value = 1.22

# Row 20
$ws.Cells.Item(20, 7).Value = 2.5
$ws.Cells.Item(20, 9).Value = 2.5
$ws.Cells.Item(20, 30).Value = 9

# Row 23
$ws.Cells.Item(23, 7).Value = 1.48
$ws.Cells.Item(23, 9).Value = 5
$ws.Cells.Item(23, 10).Value = 19
$ws.Cells.Item(23, 11).Value = 1.03
$ws.Cells.Item(23, 14).Value = 1.5
$ws.Cells.Item(23, 15).Value = 2.5
$ws.Cells.Item(23, 18).Value = 1.62
$ws.Cells.Item(23, 19).Value = 2.2
$ws.Cells.Item(23, 21).Value = 9
$ws.Cells.Item(23, 24).Value = 11
$ws.Cells.Item(23, 25).Value = 21
$ws.Cells.Item(23, 26).Value = 19
$ws.Cells.Item(23, 30).Value = 21
$ws.Cells.Item(23, 31).Value = 34
$ws.Cells.Item(23, 34).Value = 41

# Row 24
$ws.Cells.Item(24, 7).Value = 1.8
$ws.Cells.Item(24, 12).Value = 1.25
$ws.Cells.Item(24, 13).Value = 3.75
$ws.Cells.Item(24, 14).Value = 1.85
$ws.Cells.Item(24, 15).Value = 1.95

# Row 25
$ws.Cells.Item(25, 16).Value = 1.2

# Row 29
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = 3.2
$ws.Cells.Item(29, 9).Value = 3.55
$ws.Cells.Item(29, 12).Value = 1.33
$ws.Cells.Item(29, 13).Value = 2.8
$ws.Cells.Item(29, 14).Value = 1.98
$ws.Cells.Item(29, 15).Value = 1.65
$ws.Cells.Item(29, 16).Value = 1.39
$ws.Cells.Item(29, 17).Value = 2.57
$ws.Cells.Item(29, 18).Value = 1.78
$ws.Cells.Item(29, 19).Value = 1.83
$ws.Cells.Item(29, 20).Value = 6.7
$ws.Cells.Item(29, 22).Value = 8.5
$ws.Cells.Item(29, 23).Value = 18
$ws.Cells.Item(29, 24).Value = 17
$ws.Cells.Item(29, 25).Value = 30
$ws.Cells.Item(29, 26).Value = 8.75
$ws.Cells.Item(29, 27).Value = 6.2
$ws.Cells.Item(29, 28).Value = 15
$ws.Cells.Item(29, 29).Value = 75
$ws.Cells.Item(29, 30).Value = 10
$ws.Cells.Item(29, 31).Value = 19
$ws.Cells.Item(29, 32).Value = 12
$ws.Cells.Item(29, 33).Value = 50
$ws.Cells.Item(29, 34).Value = 35
$ws.Cells.Item(29, 35).Value = 40
$ws.Cells.Item(29, 36).Value = 600

Write-Output "Applied odds updates"